$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: re-apply per-cell formatting (font/fill) so every B/C/D cell carries the
#     exact style the refreshed leaderboard needs. We copy-format from a cell that already
#     has the desired look (pristine, original layout) onto every cell that needs it, one
#     style at a time. Pass order matters: B23 is both a format *source* (it is the only
#     originally-negative-style cell we read from) and, later, a format *target* (row 23
#     becomes XLRE, a positive value) -- so the "negative" pass must run before B23 gets
#     overwritten by the "positive" pass.
$xlPasteFormats = -4122

# style 10 (source B23)
$ws.Range("B23").Copy() | Out-Null
$ws.Range("B25").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B26").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B28").PasteSpecial($xlPasteFormats) | Out-Null

# style 7 (source C2)
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C18").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C19").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C21").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C25").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C26").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C27").PasteSpecial($xlPasteFormats) | Out-Null

# style 8 (source D2)
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D12").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D16").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D17").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D18").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D19").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D21").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C24").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D24").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D25").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D27").PasteSpecial($xlPasteFormats) | Out-Null

# style 9 (source C4)
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C12").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C16").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C17").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C23").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D23").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D26").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C28").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D28").PasteSpecial($xlPasteFormats) | Out-Null

# style 6 (source B2)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B12").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B16").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B17").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B18").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B19").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B21").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B23").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B24").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# --- Step 2: write the refreshed ticker list + ROC values (already sorted descending) ---
$ws.Range("A2").Value = "IBIT"
$ws.Range("B2").Value = 60.8
$ws.Range("A3").Value = "FXI"
$ws.Range("B3").Value = 37.36
$ws.Range("A4").Value = "XLF"
$ws.Range("B4").Value = 17.74
$ws.Range("A5").Value = "XLY"
$ws.Range("B5").Value = 16.68
$ws.Range("A6").Value = "MTUM"
$ws.Range("B6").Value = 15.87
$ws.Range("A7").Value = "GLD"
$ws.Range("B7").Value = 15.66
$ws.Range("A8").Value = "SLV"
$ws.Range("B8").Value = 12.83
$ws.Range("A9").Value = "SPYG"
$ws.Range("B9").Value = 12.18
$ws.Range("A10").Value = "QQQ"
$ws.Range("B10").Value = 11.11
$ws.Range("A11").Value = "SPY"
$ws.Range("B11").Value = 8.73
$ws.Range("A12").Value = "UUP"
$ws.Range("B12").Value = 7.96
$ws.Range("A13").Value = "XLI"
$ws.Range("B13").Value = 7.64
$ws.Range("A14").Value = "XLK"
$ws.Range("B14").Value = 7.6
$ws.Range("A15").Value = "USO"
$ws.Range("B15").Value = 6.37
$ws.Range("A16").Value = "XLU"
$ws.Range("B16").Value = 6.09
$ws.Range("A17").Value = "RSP"
$ws.Range("B17").Value = 5.36
$ws.Range("A18").Value = "GMF"
$ws.Range("B18").Value = 5.12
$ws.Range("A19").Value = "XLE"
$ws.Range("B19").Value = 3.81
$ws.Range("A20").Value = "SPYV"
$ws.Range("B20").Value = 3.4
$ws.Range("A21").Value = "IWO"
$ws.Range("B21").Value = 2.01
$ws.Range("A22").Value = "XLP"
$ws.Range("B22").Value = 1.94
$ws.Range("A23").Value = "XLRE"
$ws.Range("B23").Value = 1.53
$ws.Range("A24").Value = "IWN"
$ws.Range("B24").Value = 0.02
$ws.Range("A25").Value = "MOAT"
$ws.Range("B25").Value = -0.88
$ws.Range("A26").Value = "XLV"
$ws.Range("B26").Value = -2.03
$ws.Range("A27").Value = "XLB"
$ws.Range("B27").Value = -2.3
$ws.Range("A28").Value = "TLT"
$ws.Range("B28").Value = -5.07
